# -----------------------------------------------------------------------
# "Completando tabela de gerenciamento de equipe com todas as etapas"
#
# The original sheet tracked stages 1-2 of the team plan in rows 1-10
# (with a leftover blank row 11). This edit finishes the plan: it moves
# the "Criar tabelas..." deliverable up into stage 2 (row 7), turns the
# former stage-2 rows 7-9 into a new merged "2o, 3o, 4o Etapa" block
# (rows 8-10), and appends a brand-new "5o Etapa" block (rows 11-13)
# with its own deliverables. Row/column layout grows from 11 to 13 rows.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Unmerge the whole table body first - individual cells inside a
#    merged range cannot be written to until it is split back apart.
# ---------------------------------------------------------------------
$ws.Range("A1:G13").UnMerge()

# ---------------------------------------------------------------------
# 2) Grow the table: insert two new rows at the bottom (12 & 13) so the
#    old 11-row grid becomes 13 rows, matching the final stage-5 block.
# ---------------------------------------------------------------------
$ws.Range("A12:A13").EntireRow.Insert()

# ---------------------------------------------------------------------
# 3) Rewrite the cell values for rows 7-13 to their final content.
# ---------------------------------------------------------------------

# Row 7 - last row of stage "2o Etapa": deliverable becomes the KPI
# dashboard text (used to live on row 10); dates move in from old row 10.
$ws.Range("C7").Value = "Criar tabelas, gráficos ou dashboards com no mínimo 5 `nindicadores de desempenho `ne metas para o processo de negócio"
$ws.Range("E7").Value = 45189
$ws.Range("F7").Value = 45193

# Row 8 - first row of the new "2o, 3o, 4o Etapa" block.
$ws.Range("A8").Value = "2º, 3º, 4º Etapa"
$ws.Range("B8").Value = "Adeilton, Carlos,Gustavo Luiz e Pedro"
$ws.Range("C8").Value = "Programação de Funcionalidades"
$ws.Range("D8").Value = 45166
$ws.Range("E8").Value = 45188
$ws.Range("F8").Value = 45249
$ws.Range("G8").Value = 45249

# Row 9 - only the deliverable name remains (B/E/F cleared).
$ws.Range("B9").Value = $null
$ws.Range("C9").Value = "Planos de Testes de Funcionalidades e Usabilidade"
$ws.Range("E9").Value = $null
$ws.Range("F9").Value = $null

# Row 10 - only the deliverable name remains (B/E/F cleared).
$ws.Range("B10").Value = $null
$ws.Range("C10").Value = "Registros de Testes de Funcionalidades  e Usabilidade"
$ws.Range("E10").Value = $null
$ws.Range("F10").Value = $null

# Row 11 - first row of the brand-new "5o Etapa" block.
$ws.Range("A11").Value = "5º Etapa"
$ws.Range("B11").Value = "Adeilton, Carlos,Gustavo Luiz e Pedro"
$ws.Range("C11").Value = "Considerações Finais"
$ws.Range("D11").Value = 45250
$ws.Range("E11").Value = 45250
$ws.Range("F11").Value = 45264
$ws.Range("G11").Value = 45264

# Row 12 - extra "5o Etapa" deliverable.
$ws.Range("C12").Value = "Entrega de Vídeo de Apresentação Final e PDF usado na `nApresentação"

# Row 13 - extra "5o Etapa" deliverable.
$ws.Range("C13").Value = "Realização da Apresentação Final"

# Rows 7 and 12 hold the two multi-line deliverable texts; re-autofit
# them back to the sheet's standard height (entering a multi-line value
# bumps the row to a custom height otherwise).
$ws.Rows.Item(7).AutoFit()
$ws.Rows.Item(12).AutoFit()

# ---------------------------------------------------------------------
# 4) Formatting: reuse existing cell styles by copying formats from
#    cells that already carry the right look (blue "stage" cells,
#    plain centered text cells, centered date cells, etc).
# ---------------------------------------------------------------------

# B7/B8/B11 -> same "name, vertically centered" look already used by B10.
$ws.Range("B10").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# A8 -> same blue "stage" look already used by A2/A4.
$ws.Range("A4").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# E7/F7/D8/E8/F8/G8/D11/E11/F11/G11 -> centered date look already used by D2.
$ws.Range("D2").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("G11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# C7..C13 -> plain centered text look already used by C2.
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# A11 ("5o Etapa") is a brand-new combination: same centered/vertically
# centered look as the other stage cells, but with the blue fill added.
$ws.Range("A11").HorizontalAlignment = -4108
$ws.Range("A11").VerticalAlignment = -4108
$ws.Range("A11").Interior.Color = 15254943

# ---------------------------------------------------------------------
# 5) Re-merge the stage / shared-value cells for the final layout.
# ---------------------------------------------------------------------
$ws.Range("A2:A3").Merge()
$ws.Range("D2:D3").Merge()
$ws.Range("E2:E3").Merge()
$ws.Range("G2:G3").Merge()

$ws.Range("A4:A7").Merge()
$ws.Range("D4:D7").Merge()
$ws.Range("E4:E5").Merge()
$ws.Range("G4:G7").Merge()

$ws.Range("A8:A10").Merge()
$ws.Range("B8:B10").Merge()
$ws.Range("D8:D10").Merge()
$ws.Range("E8:E10").Merge()
$ws.Range("F8:F10").Merge()
$ws.Range("G8:G10").Merge()

$ws.Range("A11:A13").Merge()
$ws.Range("B11:B13").Merge()
$ws.Range("D11:D13").Merge()
$ws.Range("E11:E13").Merge()
$ws.Range("F11:F13").Merge()
$ws.Range("G11:G13").Merge()
